# add: print formulas option(-f)
# Adds a small "now / now + 1 day" example (A3, B4) using a custom
# date-time number format, auto-fits the columns that now hold data,
# and leaves the selection on the newly added cell (B4) - mirroring
# the author's manual edit captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the custom date/time format BEFORE writing the formulas so the
# engine doesn't also mint its own default datetime numFmt for us.
$ws.Range("A3").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"
$ws.Range("B4").NumberFormat = "yyyy/mm/dd\ hh:mm:ss"

$ws.Range("A3").Formula = "=NOW()"
$ws.Range("B4").Formula = "=A3+DAY(1)"

# Match the new used range / column widths to the data that was added.
$ws.Columns("A:B").AutoFit()

# Leave the selection where the author left it.
[void]$ws.Range("B4").Select()
